$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Liam James Payne"
$ws.Range("A5").Value = "Louis William Tomlinson"
$ws.Range("A6").Value = "Zayn Javadd Malik"
$ws.Range("A7").Value = "Júlia K"
